$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.963.49'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '2.526.08'
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("E4").Value = '  +0.08%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '578.85'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '151.46'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.31%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -0.05%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.161'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.77%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.26'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.45%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.352'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.34%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '29.41'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = '2.987.65'
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("D16").Value = '63.861.20'
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").Value = '2.537.22'
$ws.Range("E17").Value = '  +2.99%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '7.80'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.19%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.93'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.18%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.24'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.46%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '326.65'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.23'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  -0.13%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '10.04'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.64%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '65.44'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.71%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '658.04'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("E27").Value = '  +3.96%  '
$ws.Range("D28").Value = '2.654.75'
$ws.Range("E28").Value = '  +2.74%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.13%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.47'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.34%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.03'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.10%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.85'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("E34").Value = '  +0.04%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.53'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -1.07%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.78'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.22%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.52'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +1.42%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.83'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.52%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.370'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '151.99'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("B41").Value = 'EthereumClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '18.78'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.13%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '0.0₆0303'
$ws.Range("E44").Value = '  -5.38%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '158.13'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.54%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '15.43'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.24%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.63'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.85%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '20.84'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.61%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.617'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.35%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0517'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.88%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0228'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.06%  '
